# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - same bold/centered/bordered format as the rest of the
# header cells; copy formatting from an existing header cell so the new
# cells reuse the same style definition.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player (rows 2-48) shares the team's 2015 season record.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 100  # AD
    $ws.Cells.Item($row, 31).Value = 62   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
